# Auto-generated Excel COM-interop script to apply the cryptos.xlsx update
# Updates coin prices, volume percentages, and re-orders a few coin rows
# per the upstream data refresh (GitHub Actions cron job).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (never let Excel's type-inference
# reinterpret numeric-looking strings like "1.00" or "0.340" as numbers,
# which would silently drop the significant trailing zeros). The sheet's
# cells all use the workbook's default (unstyled) cell format, so we restore
# that after forcing the text number-format, keeping formatting untouched.
function Set-CellText {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

# Row 2
Set-CellText $ws.Range("D2") '60.687.81'
Set-CellText $ws.Range("E2") '  -1.59%  '
# Row 3
Set-CellText $ws.Range("D3") '2.348.87'
Set-CellText $ws.Range("E3") '  -5.26%  '
# Row 4
Set-CellText $ws.Range("E4") '  +0.14%  '
# Row 5
Set-CellText $ws.Range("D5") '542.34'
Set-CellText $ws.Range("E5") '  -2.53%  '
# Row 6
Set-CellText $ws.Range("D6") '138.53'
Set-CellText $ws.Range("E6") '  -6.02%  '
# Row 7
Set-CellText $ws.Range("E7") '  +0.05%  '
# Row 8
Set-CellText $ws.Range("D8") '0.517'
Set-CellText $ws.Range("E8") '  -14.00%  '
# Row 9
Set-CellText $ws.Range("D9") '2.348.68'
Set-CellText $ws.Range("E9") '  -5.20%  '
# Row 10
Set-CellText $ws.Range("D10") '0.104'
Set-CellText $ws.Range("E10") '  -3.92%  '
# Row 12
Set-CellText $ws.Range("D12") '5.22'
Set-CellText $ws.Range("E12") '  -4.93%  '
# Row 13
Set-CellText $ws.Range("D13") '0.340'
Set-CellText $ws.Range("E13") '  -5.12%  '
# Row 14
Set-CellText $ws.Range("D14") '24.75'
Set-CellText $ws.Range("E14") '  -6.72%  '
# Row 15
Set-CellText $ws.Range("D15") '2.781.94'
Set-CellText $ws.Range("E15") '  -4.94%  '
# Row 16
Set-CellText $ws.Range("B16") 'WrappedBTC'
Set-CellText $ws.Range("C16") 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-CellText $ws.Range("D16") '60.873.38'
Set-CellText $ws.Range("E16") '  -1.18%  '
# Row 17
Set-CellText $ws.Range("B17") 'ShibaInu'
Set-CellText $ws.Range("C17") 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-CellText $ws.Range("D17") '0.0000161'
Set-CellText $ws.Range("E17") '  -4.37%  '
# Row 18
Set-CellText $ws.Range("D18") '2.361.74'
Set-CellText $ws.Range("E18") '  -4.95%  '
# Row 19
Set-CellText $ws.Range("D19") '10.55'
Set-CellText $ws.Range("E19") '  -5.82%  '
# Row 20
Set-CellText $ws.Range("D20") '4.06'
Set-CellText $ws.Range("E20") '  -4.15%  '
# Row 21
Set-CellText $ws.Range("D21") '313.58'
Set-CellText $ws.Range("E21") '  -2.63%  '
# Row 22
Set-CellText $ws.Range("D22") '6.56'
Set-CellText $ws.Range("E22") '  -8.84%  '
# Row 23
Set-CellText $ws.Range("D23") '0.998'
Set-CellText $ws.Range("E23") '  -0.24%  '
# Row 24
Set-CellText $ws.Range("D24") '1.84'
Set-CellText $ws.Range("E24") '  -2.96%  '
# Row 25
Set-CellText $ws.Range("D25") '63.12'
Set-CellText $ws.Range("E25") '  -2.03%  '
# Row 26
Set-CellText $ws.Range("D26") '8.03'
Set-CellText $ws.Range("E26") '  +3.34%  '
# Row 27
Set-CellText $ws.Range("E27") '  +0.23%  '
# Row 28
Set-CellText $ws.Range("D28") '2.480.30'
Set-CellText $ws.Range("E28") '  -4.93%  '
# Row 29
Set-CellText $ws.Range("D29") '0.0₃0899'
Set-CellText $ws.Range("E29") '  -10.40%  '
# Row 30
Set-CellText $ws.Range("D30") '511.77'
Set-CellText $ws.Range("E30") '  -9.16%  '
# Row 31
Set-CellText $ws.Range("D31") '7.92'
Set-CellText $ws.Range("E31") '  -5.39%  '
# Row 32
Set-CellText $ws.Range("D32") '1.38'
Set-CellText $ws.Range("E32") '  -9.22%  '
# Row 33
Set-CellText $ws.Range("D33") '0.144'
Set-CellText $ws.Range("E33") '  -4.64%  '
# Row 34
Set-CellText $ws.Range("D34") '1.81'
Set-CellText $ws.Range("E34") '  -6.68%  '
# Row 35
Set-CellText $ws.Range("D35") '1.53'
Set-CellText $ws.Range("E35") '  -5.16%  '
# Row 36
Set-CellText $ws.Range("D36") '1.00'
Set-CellText $ws.Range("E36") '  +0.10%  '
# Row 37
Set-CellText $ws.Range("D37") '4.57'
Set-CellText $ws.Range("E37") '  -7.70%  '
# Row 38
Set-CellText $ws.Range("B38") 'PolygonEcosystemToken'
Set-CellText $ws.Range("C38") 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-CellText $ws.Range("D38") '0.370'
Set-CellText $ws.Range("E38") '  -3.82%  '
# Row 39
Set-CellText $ws.Range("B39") 'RenderToken'
Set-CellText $ws.Range("C39") 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-CellText $ws.Range("D39") '5.26'
Set-CellText $ws.Range("E39") '  -11.84%  '
# Row 40
Set-CellText $ws.Range("B40") 'EthereumClassic'
Set-CellText $ws.Range("C40") 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-CellText $ws.Range("D40") '18.12'
Set-CellText $ws.Range("E40") '  -2.70%  '
# Row 41
Set-CellText $ws.Range("B41") 'Stacks'
Set-CellText $ws.Range("C41") 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-CellText $ws.Range("D41") '1.80'
Set-CellText $ws.Range("E41") '  +1.25%  '
# Row 42
Set-CellText $ws.Range("D42") '139.26'
Set-CellText $ws.Range("E42") '  -4.20%  '
# Row 43
Set-CellText $ws.Range("E43") '  +0.01%  '
# Row 44
Set-CellText $ws.Range("D44") '40.09'
Set-CellText $ws.Range("E44") '  -1.35%  '
# Row 45
Set-CellText $ws.Range("B45") 'dogwifhat'
Set-CellText $ws.Range("C45") 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-CellText $ws.Range("D45") '2.10'
Set-CellText $ws.Range("E45") '  -14.71%  '
# Row 46
Set-CellText $ws.Range("B46") 'Filecoin'
Set-CellText $ws.Range("C46") 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-CellText $ws.Range("D46") '3.55'
Set-CellText $ws.Range("E46") '  -2.82%  '
# Row 47
Set-CellText $ws.Range("B47") 'Aave'
Set-CellText $ws.Range("C47") 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-CellText $ws.Range("D47") '137.78'
Set-CellText $ws.Range("E47") '  -7.29%  '
# Row 48
Set-CellText $ws.Range("D48") '0.0509'
Set-CellText $ws.Range("E48") '  -6.59%  '
# Row 49
Set-CellText $ws.Range("D49") '19.49'
Set-CellText $ws.Range("E49") '  -12.12%  '
# Row 50
Set-CellText $ws.Range("D50") '0.565'
Set-CellText $ws.Range("E50") '  -5.64%  '
# Row 51
Set-CellText $ws.Range("D51") '0.0894'
Set-CellText $ws.Range("E51") '  -5.38%  '

"Applied cryptos.xlsx update: " + 49 + " rows touched"
